$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.657.19"
$ws.Range("E2").Value = "  +3.69%  "
$ws.Range("D3").Value = "3.503.77"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'591.66"
$ws.Range("E5").Value = "  +3.11%  "
$ws.Range("D6").Value = "'169.69"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  +8.04%  "
$ws.Range("D9").Value = "3.501.07"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "'0.128"
$ws.Range("E11").Value = "  +6.08%  "
$ws.Range("D12").Value = "'0.440"
$ws.Range("E12").Value = "  +3.20%  "
$ws.Range("D13").Value = "4.109.20"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "66.664.05"
$ws.Range("E17").Value = "  +3.61%  "
$ws.Range("D18").Value = "3.482.64"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").Value = "'14.16"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").Value = "'392.85"
$ws.Range("E21").Value = "  +3.46%  "
$ws.Range("D22").Value = "'8.00"
$ws.Range("E22").Value = "  +1.96%  "
$ws.Range("D23").Value = "'73.24"
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("E26").Value = "  +4.83%  "
$ws.Range("D27").Value = "'10.24"
$ws.Range("E27").Value = "  +6.54%  "
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "'6.36"
$ws.Range("E30").Value = "  +3.99%  "
$ws.Range("E31").Value = "  +4.18%  "
$ws.Range("E32").Value = "  +2.81%  "
$ws.Range("D33").Value = "'23.65"
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +7.55%  "
$ws.Range("D37").Value = "'162.19"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").Value = "'0.886"
$ws.Range("E38").Value = "  +2.63%  "
$ws.Range("E39").Value = "  +4.32%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'6.83"
$ws.Range("E40").Value = "  +5.37%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'27.69"
$ws.Range("E41").Value = "  +4.07%  "
$ws.Range("E42").Value = "  +5.57%  "
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").Value = "2.794.75"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "'43.28"
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("D47").Value = "'2.53"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").Value = "'0.0312"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").Value = "'351.65"
$ws.Range("E49").Value = "  +4.95%  "
$ws.Range("E50").Value = "  +4.15%  "
$ws.Range("D51").Value = "'33.71"
